# Adds a new "data contact us" worksheet (after "data laptop") with the
# email/subject test-case rows used by the "assertion" test, matching
# the QA test-data workbook used in the "Final-Project-QA-WEB" repo.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet ("data laptop") so it
# lands at the end of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "data contact us"

# Header row
$ws.Range("A1").Value = "email"
$ws.Range("B1").Value = "subject"

# Row 2 - first contact-us test record
$ws.Range("A2").Value = "yudhatesting@gmail.com"
$ws.Range("B2").Value = "this product so cool!"

# Row 3 - second contact-us test record
$ws.Range("A3").Value = "yudhatesting2@gmail.com"
$ws.Range("B3").Value = "I want to buy this product"

# Turn the two email addresses into live mailto: hyperlinks, then re-apply
# the "Hyperlink" cell style (Hyperlinks.Add already does this, this just
# keeps the styling explicit/deterministic).
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:yudhatesting@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:yudhatesting2@gmail.com")
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("A3").Style = "Hyperlink"

# Column widths to fit the email / subject text
$ws.Columns.Item(1).ColumnWidth = 20.6080729166667
$ws.Columns.Item(2).ColumnWidth = 24.4986979166667

# Make this the active sheet/selection (tabSelected + activeCell = B3)
$ws.Range("B3").Select()
$ws.Activate()
